$wb = $excel.ActiveWorkbook

# --- Sheet "BpTPEU" edits ---
$wsData = $wb.Worksheets.Item("BpTPEU")

# B1: "One Quadrillion BTU" -> "PJ", also highlight with yellow fill (keep right alignment)
$wsData.Range("B1").Value = "PJ"
$wsData.Range("B1").Interior.Color = 65535
$wsData.Range("B1").HorizontalAlignment = -4152   # xlRight

# B2: formula changes from =10^15 to =9.478*10^11
$wsData.Range("B2").Formula = "=9.478*10^11"

# --- Sheet "About" edits ---
$wsAbout = $wb.Worksheets.Item("About")

# B11: update the Brazilian adaptation note text
$wsAbout.Range("B11").Value = "We assumed that the Total Primary Energy output is PJ for Brazil. "

# Extend the yellow highlight fill one more cell to the right (G11)
$wsAbout.Range("G11").Interior.Color = 65535

# --- Restore / update the on-screen selections to match the saved view state ---
# Touch BpTPEU's selection first, then re-activate About last so "About" stays
# the selected/visible tab (as it was before the edit).
$wsData.Activate()
$wsData.Range("B22").Select()

$wsAbout.Activate()
$wsAbout.Range("H16").Select()
